# modular-operations-scaling.xlsx
#
# "Avoid vector operation when value at pivot = 0"
#
# The underlying experiment was re-run with a fix that skips a vector
# operation whenever the pivot value is zero. That changed the measured
# operation counts (column B) and the fitted exponent used in C2 (which in
# turn drives the extrapolated "N^C" column C). The scratch Log(N)/Log(Ops)
# helper columns (E/F) that were used to eyeball-fit the exponent are no
# longer needed and are removed, along with the now-stale "chart data
# range" defined names Excel had created for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the hidden "_xlchart.v1.*" defined names. These were generated
#    by Excel to remember the chart's filtered source ranges and are no
#    longer valid once the helper columns go away.
# ---------------------------------------------------------------------
$namesToDelete = @()
foreach ($n in $wb.Names) {
    if ($n.Name -like "_xlchart*") {
        $namesToDelete += $n.Name
    }
}
foreach ($nm in $namesToDelete) {
    $wb.Names.Item($nm).Delete()
}

# ---------------------------------------------------------------------
# 2. Update the fitted exponent in C2 (3.9 -> 2.15) and give it a more
#    precise display format (three decimal places instead of one).
# ---------------------------------------------------------------------
$ws.Range("C2").Value = 2.15
$ws.Range("C2").NumberFormat = "0.000"

# ---------------------------------------------------------------------
# 3. Rename the C-column header from "N^3.9" to the more generic "N^C"
#    (it now reflects whatever exponent happens to be in C2).
# ---------------------------------------------------------------------
$ws.Range("C3").Value = "N^C"

# ---------------------------------------------------------------------
# 4. Replace the measured operation counts (column B, rows 4-16) with the
#    re-measured values from the corrected run.
# ---------------------------------------------------------------------
$newOps = @(1041, 1719, 2565, 3579, 4761, 7629, 11169, 15381, 20265, 32049, 46521, 63681, 83529)
for ($i = 0; $i -lt $newOps.Length; $i++) {
    $row = 4 + $i
    $ws.Range("B$row").Value = $newOps[$i]
}
# Column C (the extrapolated N^C fit) keeps its existing formula, which
# recalculates automatically against the new B values / C2 exponent.

# ---------------------------------------------------------------------
# 5. Remove the scratch "Log N" / "Log Ops" helper columns (E & F) that
#    were only used to manually eyeball the old exponent.
# ---------------------------------------------------------------------
$ws.Range("E1:F16").Clear()

# ---------------------------------------------------------------------
# 6. Tidy up the view: drop the frozen/scrolled top-left cell and move
#    the active selection onto the updated C4 fit cell.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("C4").Select()
